$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40..172 down to 41..173
$ws.Rows("40:40").Insert()

# Populate the newly inserted row 40 with the new record
$ws.Range("A40").Value = 5
$ws.Range("B40").Value = "Macroferia Regional de Talca"
$ws.Range("C40").Value = "Maule"
$ws.Range("D40").Value = 44481
$ws.Range("E40").Value = 7
$ws.Range("F40").Value = 100112003
$ws.Range("G40").Value = "Ajo"
$ws.Range("H40").Value = "Chino"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 300
$ws.Range("K40").Value = 15000
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = 15000
$ws.Range("N40").Value = "$/caja 10 kilos"
$ws.Range("O40").Value = "China"
$ws.Range("P40").Value = 1500
$ws.Range("Q40").Value = 10
$ws.Range("R40").Value = "Hortaliza"
